$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.488.24"
$ws.Range("E2").Value = "  +0.96%  "

$ws.Range("D3").Value = "'1.876.54"
$ws.Range("E3").Value = "  +1.28%  "

$ws.Range("D4").Value = "'0.9990"
$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").Value = "'0.7160"
$ws.Range("E5").Value = "  +2.33%  "

$ws.Range("D6").Value = "'241.56"
$ws.Range("E6").Value = "  +1.67%  "

$ws.Range("D7").Value = "'1.0000"
$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("D8").Value = "'0.07954"
$ws.Range("E8").Value = "  +0.83%  "

$ws.Range("D9").Value = "'0.3107"
$ws.Range("E9").Value = "  +3.03%  "

$ws.Range("D10").Value = "'25.35"
$ws.Range("E10").Value = "  +5.92%  "

$ws.Range("D11").Value = "'0.08288"

$ws.Range("D12").Value = "'0.7301"
$ws.Range("E12").Value = "  +3.35%  "

$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "'1.878.88"
$ws.Range("E13").Value = "  +1.21%  "

$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "'5.284"
$ws.Range("E14").Value = "  +1.84%  "

$ws.Range("D15").Value = "'91.22"
$ws.Range("E15").Value = "  +2.05%  "

$ws.Range("D16").Value = "'29.466.57"
$ws.Range("E16").Value = "  +0.92%  "

$ws.Range("D17").Value = "'5.934"
$ws.Range("E17").Value = "  +2.18%  "

$ws.Range("D18").Value = "'245.60"
$ws.Range("E18").Value = "  +4.25%  "

$ws.Range("D19").Value = "'0.000007884"
$ws.Range("E19").Value = "  +0.69%  "

$ws.Range("D20").Value = "'13.36"

$ws.Range("D21").Value = "'2.117.58"
$ws.Range("E21").Value = "  +1.25%  "

$ws.Range("D22").Value = "'0.9995"
$ws.Range("E22").Value = "  -0.05%  "

$ws.Range("D23").Value = "'7.984"
$ws.Range("E23").Value = "  +6.51%  "

$ws.Range("D24").Value = "'1.001"
$ws.Range("E24").Value = "  +0.08%  "

$ws.Range("D25").Value = "'0.1608"
$ws.Range("E25").Value = "  +13.50%  "

$ws.Range("D26").Value = "'163.73"
$ws.Range("E26").Value = "  +0.73%  "

$ws.Range("D27").Value = "'9.054"
$ws.Range("E27").Value = "  +1.91%  "

$ws.Range("D28").Value = "'18.33"
$ws.Range("E28").Value = "  +1.68%  "

$ws.Range("D29").Value = "'1.359"
$ws.Range("E29").Value = "  -2.98%  "

$ws.Range("D30").Value = "'1.490"
$ws.Range("E30").Value = "  +1.35%  "

$ws.Range("D31").Value = "'4.391"
$ws.Range("E31").Value = "  +2.09%  "

$ws.Range("D32").Value = "'4.116"
$ws.Range("E32").Value = "  +2.65%  "

$ws.Range("D33").Value = "'0.05267"
$ws.Range("E33").Value = "  +2.38%  "

$ws.Range("D34").Value = "'1.957"
$ws.Range("E34").Value = "  +2.12%  "

$ws.Range("E35").Value = "  +2.70%  "

$ws.Range("D36").Value = "'0.7272"
$ws.Range("E36").Value = "  +2.98%  "

$ws.Range("D37").Value = "'2.675"
$ws.Range("E37").Value = "  -0.07%  "

$ws.Range("E38").Value = "  +1.36%  "

$ws.Range("D39").Value = "'1.221.64"

$ws.Range("D40").Value = "'2.703"
$ws.Range("E40").Value = "  -0.03%  "

$ws.Range("D41").Value = "'0.9112"
$ws.Range("E41").Value = "  -1.03%  "

$ws.Range("D42").Value = "'73.96"
$ws.Range("E42").Value = "  +5.58%  "

$ws.Range("D43").Value = "'6.123"
$ws.Range("E43").Value = "  +2.99%  "

$ws.Range("D44").Value = "'1.000"
$ws.Range("E44").Value = "  +0.01%  "

$ws.Range("D45").Value = "'102.22"
$ws.Range("E45").Value = "  -0.71%  "

$ws.Range("D46").Value = "'0.5287"
$ws.Range("E46").Value = "  -0.04%  "

$ws.Range("D47").Value = "'2.012.24"
$ws.Range("E47").Value = "  +0.81%  "

$ws.Range("D48").Value = "'1.799"
$ws.Range("E48").Value = "  +3.61%  "

$ws.Range("D49").Value = "'2.927"
$ws.Range("E49").Value = "  +9.18%  "

$ws.Range("E50").Value = "  +1.64%  "

$ws.Range("D51").Value = "'9.346"
$ws.Range("E51").Value = "  +2.04%  "
